# VACCINECERT-1633 Fixed CSV upload templates
#
# The "vaccination" CSV-upload sample/template had a stale example
# sampleDate of 2021-11-01 (serial 44501) in the demo data row. Update it
# to 2021-11-16 (serial 44516) to match the fixed template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the sample data row's sampleDate value.
$ws.Range("E2").Value = 44516

# Leave the selection where it was when the fixed template was saved.
$ws.Range("E10").Select()
